$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.960.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "'1.574.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.89%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").Value = "'299.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.82%  "

$ws.Range("D7").Value = "'0.3752"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.3553"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.99%  "

$ws.Range("D9").Value = "'49.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.37%  "

$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("D11").Value = "'1.212"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.64%  "

$ws.Range("D12").Value = "'0.07953"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.06%  "

$ws.Range("D13").Value = "'21.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.75%  "

$ws.Range("D14").Value = "'6.384"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.24%  "

$ws.Range("D15").Value = "'7.280"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.77%  "

$ws.Range("D16").Value = "'0.00001222"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.08%  "

$ws.Range("D17").Value = "'1.575.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.90%  "

$ws.Range("D18").Value = "'91.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.66%  "

$ws.Range("D19").Value = "'0.06735"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("D20").Value = "'17.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.46%  "

$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").Value = "'6.337"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.13%  "

$ws.Range("D23").Value = "'22.947.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("D24").Value = "'12.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.60%  "

$ws.Range("D25").Value = "'2.373"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("D26").Value = "'2.808"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.51%  "

$ws.Range("D27").Value = "'20.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.18%  "

$ws.Range("D28").Value = "'148.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.40%  "

$ws.Range("D29").Value = "'5.166"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.68%  "

$ws.Range("D30").Value = "'131.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.77%  "

$ws.Range("D31").Value = "'2.331"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.44%  "

$ws.Range("D32").Value = "'6.534"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.24%  "

$ws.Range("D33").Value = "'1.749.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.86%  "

$ws.Range("D34").Value = "'0.9289"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.89%  "

$ws.Range("E35").Value = "  -4.75%  "

$ws.Range("D36").Value = "'0.08746"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.95%  "

$ws.Range("D37").Value = "'9.913"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.16%  "

$ws.Range("D38").Value = "'0.02627"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.30%  "

$ws.Range("D39").Value = "'0.2454"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.97%  "

$ws.Range("D40").Value = "'5.953"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.62%  "

$ws.Range("D41").Value = "'1.342"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.45%  "

$ws.Range("D42").Value = "'0.6845"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.65%  "

$ws.Range("D43").Value = "'11.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.15%  "

$ws.Range("D44").Value = "'14.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.85%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").Value = "'0.6304"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.90%  "

$ws.Range("D47").Value = "'3.960"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("D48").Value = "'2.240"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.68%  "

$ws.Range("D49").Value = "'0.07839"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.81%  "

$ws.Range("D50").Value = "'129.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.27%  "

$ws.Range("D51").Value = "'1.179"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.50%  "
